$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-31 Saturday" "2026-02-01 Sunday"

Replace-Text "259×9=" "811×5="
Replace-Text "381×7=" "275×9="
Replace-Text "710×8=" "748×7="
Replace-Text "596×8=" "187×7="
Replace-Text "621×6=" "458×5="

Replace-Text "660×9=" "326×9="
Replace-Text "606×5=" "193×5="
Replace-Text "638×5=" "463×5="
Replace-Text "559×8=" "355×7="
Replace-Text "954×2=" "922×7="

Replace-Text "361×6=" "466×3="
Replace-Text "458×9=" "439×9="
Replace-Text "782×9=" "435×8="
Replace-Text "519×7=" "218×3="
Replace-Text "319×6=" "273×9="

Replace-Text "770×8=" "623×2="
Replace-Text "229×6=" "897×5="
Replace-Text "331×5=" "777×2="
Replace-Text "834×6=" "841×3="
Replace-Text "328×8=" "588×8="

Replace-Text "175×4=" "729×8="
Replace-Text "970×2=" "458×8="
Replace-Text "154×4=" "911×8="
Replace-Text "428×9=" "631×5="
Replace-Text "722×7=" "766×8="

Write-Output "Done"
